$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText {
    param($cellRef, $text)
    # Force the new value to be stored as plain text (matching the original
    # inline-string cell) rather than letting Excel auto-detect it as a date
    # serial number. Reset the style back to Normal afterwards so no stray
    # NumberFormat/style is left behind on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 3: date format change + D/G count updates
Set-DateText "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: date format change only
Set-DateText "A4" "01-08-2022"

# Row 5: date format change + D/E/H count updates
Set-DateText "A5" "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6: date format change only
Set-DateText "A6" "08-08-2022"

# Row 7: date format change only
Set-DateText "A7" "11-08-2022"

# Row 8: date format change only
Set-DateText "A8" "15-08-2022"

# Row 9: date format change only
Set-DateText "A9" "18-08-2022"

# Row 10: date format change only
Set-DateText "A10" "22-08-2022"

# Row 11: date format change only
Set-DateText "A11" "25-08-2022"

# Row 12: date format change only
Set-DateText "A12" "29-08-2022"

# Row 13: date format change only
Set-DateText "A13" "01-09-2022"

# Row 14: date format change only
Set-DateText "A14" "05-09-2022"

# Row 15: date format change only
Set-DateText "A15" "08-09-2022"

# Row 16: date format change only
Set-DateText "A16" "12-09-2022"

# Row 17: date format change only
Set-DateText "A17" "15-09-2022"

# Row 18: date format change only
Set-DateText "A18" "19-09-2022"

# Row 19: date format change only
Set-DateText "A19" "22-09-2022"

# Row 20: date format change only
Set-DateText "A20" "26-09-2022"

# Row 21: date format change only
Set-DateText "A21" "29-09-2022"
